# Update "想去人数" (want-to-go count) figures to the latest scrape values.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 953
$wsExhibit.Range("F3").Value = 1851
$wsExhibit.Range("F4").Value = 415

# Sheet "全部类型" (All types) contains the same three exhibition rows, offset by the
# two "演出" rows that precede them.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 953
$wsAll.Range("F5").Value = 1851
$wsAll.Range("F6").Value = 415
